$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the style of the existing header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values per row (2-22)
$saveValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
